# Update countries & provincias Spain
# Refresh the COVID-19 country table ('Pais' sheet) with the newer data pull:
#  - several countries' case/death/recovery counts increase
#  - the table (sorted descending by 'Casos totales', column B) is re-sorted,
#    which re-orders a handful of adjacent rows (e.g. Bolivia/Venezuela,
#    Malaui/Sahara Occidental, Anguila/Islas Virgenes Britanicas, Islas Malvinas)
#  - the 'last updated' timestamp footer moves from 20:22 to 20:52

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 20:52"

# Full refreshed + re-sorted country table, written in one shot as a 2-D array
# (columns: Pais, Casos totales, Nuevos casos, Casos activos, Recuperados,
#  Casos criticos, Muertes hoy, Muertes) starting at A4.
$data = New-Object 'object[,]' 211,8
$data[0,0] = 'Estados Unidos'; $data[0,1] = 328662; $data[0,2] = 17305; $data[0,3] = 16700; $data[0,4] = 302597; $data[0,5] = 8542; $data[0,6] = 913; $data[0,7] = 9365
$data[1,0] = 'España'; $data[1,1] = 130759; $data[1,2] = 4591; $data[1,3] = 38080; $data[1,4] = 80261; $data[1,5] = 6861; $data[1,6] = 471; $data[1,7] = 12418
$data[2,0] = 'Italia'; $data[2,1] = 128948; $data[2,2] = 4316; $data[2,3] = 21815; $data[2,4] = 91246; $data[2,5] = 3977; $data[2,6] = 525; $data[2,7] = 15887
$data[3,0] = 'Alemania'; $data[3,1] = 99964; $data[3,2] = 3872; $data[3,3] = 28700; $data[3,4] = 69691; $data[3,5] = 3936; $data[3,6] = 129; $data[3,7] = 1573
$data[4,0] = 'Francia'; $data[4,1] = 89953; $data[4,2] = 0; $data[4,3] = 15438; $data[4,4] = 66955; $data[4,5] = 6838; $data[4,6] = 0; $data[4,7] = 7560
$data[5,0] = 'China'; $data[5,1] = 81669; $data[5,2] = 30; $data[5,3] = 76964; $data[5,4] = 1376; $data[5,5] = 295; $data[5,6] = 3; $data[5,7] = 3329
$data[6,0] = 'Iran'; $data[6,1] = 58226; $data[6,2] = 2483; $data[6,3] = 19736; $data[6,4] = 34887; $data[6,5] = 4103; $data[6,6] = 151; $data[6,7] = 3603
$data[7,0] = 'Reino Unido'; $data[7,1] = 47806; $data[7,2] = 5903; $data[7,3] = 135; $data[7,4] = 42737; $data[7,5] = 1559; $data[7,6] = 621; $data[7,7] = 4934
$data[8,0] = 'Turquia'; $data[8,1] = 27069; $data[8,2] = 3135; $data[8,3] = 1042; $data[8,4] = 25453; $data[8,5] = 1381; $data[8,6] = 73; $data[8,7] = 574
$data[9,0] = 'Suiza'; $data[9,1] = 21100; $data[9,2] = 595; $data[9,3] = 6415; $data[9,4] = 14000; $data[9,5] = 391; $data[9,6] = 19; $data[9,7] = 685
$data[10,0] = 'Belgica'; $data[10,1] = 19691; $data[10,2] = 1260; $data[10,3] = 3751; $data[10,4] = 14493; $data[10,5] = 1261; $data[10,6] = 164; $data[10,7] = 1447
$data[11,0] = 'Paises Bajos'; $data[11,1] = 17851; $data[11,2] = 1224; $data[11,3] = 250; $data[11,4] = 15835; $data[11,5] = 1360; $data[11,6] = 115; $data[11,7] = 1766
$data[12,0] = 'Canada'; $data[12,1] = 15416; $data[12,2] = 1504; $data[12,3] = 2613; $data[12,4] = 12526; $data[12,5] = 426; $data[12,6] = 46; $data[12,7] = 277
$data[13,0] = 'Austria'; $data[13,1] = 12051; $data[13,2] = 270; $data[13,3] = 2998; $data[13,4] = 8849; $data[13,5] = 244; $data[13,6] = 18; $data[13,7] = 204
$data[14,0] = 'Portugal'; $data[14,1] = 11278; $data[14,2] = 754; $data[14,3] = 75; $data[14,4] = 10908; $data[14,5] = 267; $data[14,6] = 29; $data[14,7] = 295
$data[15,0] = 'Brasil'; $data[15,1] = 10568; $data[15,2] = 208; $data[15,3] = 127; $data[15,4] = 9986; $data[15,5] = 296; $data[15,6] = 10; $data[15,7] = 455
$data[16,0] = 'Corea del Sur'; $data[16,1] = 10237; $data[16,2] = 81; $data[16,3] = 6463; $data[16,4] = 3591; $data[16,5] = 55; $data[16,6] = 6; $data[16,7] = 183
$data[17,0] = 'Israel'; $data[17,1] = 8430; $data[17,2] = 579; $data[17,3] = 477; $data[17,4] = 7904; $data[17,5] = 139; $data[17,6] = 5; $data[17,7] = 49
$data[18,0] = 'Suecia'; $data[18,1] = 6830; $data[18,2] = 387; $data[18,3] = 205; $data[18,4] = 6224; $data[18,5] = 541; $data[18,6] = 28; $data[18,7] = 401
$data[19,0] = 'Noruega'; $data[19,1] = 5687; $data[19,2] = 137; $data[19,3] = 32; $data[19,4] = 5584; $data[19,5] = 89; $data[19,6] = 9; $data[19,7] = 71
$data[20,0] = 'Australia'; $data[20,1] = 5687; $data[20,2] = 137; $data[20,3] = 2315; $data[20,4] = 3337; $data[20,5] = 91; $data[20,6] = 5; $data[20,7] = 35
$data[21,0] = 'Rusia'; $data[21,1] = 5389; $data[21,2] = 658; $data[21,3] = 355; $data[21,4] = 4989; $data[21,5] = 8; $data[21,6] = 2; $data[21,7] = 45
$data[22,0] = 'Irlanda'; $data[22,1] = 4994; $data[22,2] = 390; $data[22,3] = 25; $data[22,4] = 4811; $data[22,5] = 165; $data[22,6] = 21; $data[22,7] = 158
$data[23,0] = 'Chequia'; $data[23,1] = 4543; $data[23,2] = 71; $data[23,3] = 96; $data[23,4] = 4380; $data[23,5] = 86; $data[23,6] = 8; $data[23,7] = 67
$data[24,0] = 'Chile'; $data[24,1] = 4471; $data[24,2] = 310; $data[24,3] = 618; $data[24,4] = 3819; $data[24,5] = 307; $data[24,6] = 7; $data[24,7] = 34
$data[25,0] = 'Dinamarca'; $data[25,1] = 4369; $data[25,2] = 292; $data[25,3] = 1327; $data[25,4] = 2863; $data[25,5] = 142; $data[25,6] = 18; $data[25,7] = 179
$data[26,0] = 'Polonia'; $data[26,1] = 4102; $data[26,2] = 475; $data[26,3] = 134; $data[26,4] = 3874; $data[26,5] = 50; $data[26,6] = 15; $data[26,7] = 94
$data[27,0] = 'Rumania'; $data[27,1] = 3864; $data[27,2] = 251; $data[27,3] = 374; $data[27,4] = 3339; $data[27,5] = 141; $data[27,6] = 5; $data[27,7] = 151
$data[28,0] = 'Malasia'; $data[28,1] = 3662; $data[28,2] = 179; $data[28,3] = 1005; $data[28,4] = 2596; $data[28,5] = 99; $data[28,6] = 4; $data[28,7] = 61
$data[29,0] = 'Ecuador'; $data[29,1] = 3646; $data[29,2] = 181; $data[29,3] = 100; $data[29,4] = 3366; $data[29,5] = 100; $data[29,6] = 8; $data[29,7] = 180
$data[30,0] = 'India'; $data[30,1] = 3588; $data[30,2] = 0; $data[30,3] = 229; $data[30,4] = 3260; $data[30,5] = 0; $data[30,6] = 0; $data[30,7] = 99
$data[31,0] = 'Filipinas'; $data[31,1] = 3246; $data[31,2] = 152; $data[31,3] = 64; $data[31,4] = 3030; $data[31,5] = 1; $data[31,6] = 8; $data[31,7] = 152
$data[32,0] = 'Japon'; $data[32,1] = 3139; $data[32,2] = 0; $data[32,3] = 514; $data[32,4] = 2548; $data[32,5] = 64; $data[32,6] = 0; $data[32,7] = 77
$data[33,0] = 'Pakistan'; $data[33,1] = 3123; $data[33,2] = 305; $data[33,3] = 208; $data[33,4] = 2870; $data[33,5] = 18; $data[33,6] = 4; $data[33,7] = 45
$data[34,0] = 'Luxemburgo'; $data[34,1] = 2804; $data[34,2] = 75; $data[34,3] = 500; $data[34,4] = 2268; $data[34,5] = 33; $data[34,6] = 5; $data[34,7] = 36
$data[35,0] = 'Arabia Saudita'; $data[35,1] = 2385; $data[35,2] = 206; $data[35,3] = 488; $data[35,4] = 1863; $data[35,5] = 41; $data[35,6] = 5; $data[35,7] = 34
$data[36,0] = 'Peru'; $data[36,1] = 2281; $data[36,2] = 535; $data[36,3] = 989; $data[36,4] = 1209; $data[36,5] = 81; $data[36,6] = 10; $data[36,7] = 83
$data[37,0] = 'Indonesia'; $data[37,1] = 2273; $data[37,2] = 181; $data[37,3] = 164; $data[37,4] = 1911; $data[37,5] = 0; $data[37,6] = 7; $data[37,7] = 198
$data[38,0] = 'Tailandia'; $data[38,1] = 2169; $data[38,2] = 102; $data[38,3] = 793; $data[38,4] = 1353; $data[38,5] = 23; $data[38,6] = 3; $data[38,7] = 23
$data[39,0] = 'Finlandia'; $data[39,1] = 1927; $data[39,2] = 45; $data[39,3] = 300; $data[39,4] = 1599; $data[39,5] = 73; $data[39,6] = 3; $data[39,7] = 28
$data[40,0] = 'Serbia'; $data[40,1] = 1908; $data[40,2] = 284; $data[40,3] = 54; $data[40,4] = 1803; $data[40,5] = 98; $data[40,6] = 7; $data[40,7] = 51
$data[41,0] = 'Mexico'; $data[41,1] = 1890; $data[41,2] = 202; $data[41,3] = 633; $data[41,4] = 1178; $data[41,5] = 1; $data[41,6] = 19; $data[41,7] = 79
$data[42,0] = 'Panama'; $data[42,1] = 1801; $data[42,2] = 0; $data[42,3] = 13; $data[42,4] = 1742; $data[42,5] = 75; $data[42,6] = 0; $data[42,7] = 46
$data[43,0] = 'Emiratos Arabes Unidos'; $data[43,1] = 1798; $data[43,2] = 293; $data[43,3] = 144; $data[43,4] = 1644; $data[43,5] = 1; $data[43,6] = 0; $data[43,7] = 10
$data[44,0] = 'Republica Dominicana'; $data[44,1] = 1745; $data[44,2] = 167; $data[44,3] = 17; $data[44,4] = 1646; $data[44,5] = 147; $data[44,6] = 5; $data[44,7] = 82
$data[45,0] = 'Grecia'; $data[45,1] = 1735; $data[45,2] = 62; $data[45,3] = 78; $data[45,4] = 1584; $data[45,5] = 93; $data[45,6] = 5; $data[45,7] = 73
$data[46,0] = 'Catar'; $data[46,1] = 1604; $data[46,2] = 279; $data[46,3] = 123; $data[46,4] = 1477; $data[46,5] = 37; $data[46,6] = 1; $data[46,7] = 4
$data[47,0] = 'Sudafrica'; $data[47,1] = 1585; $data[47,2] = 0; $data[47,3] = 95; $data[47,4] = 1481; $data[47,5] = 7; $data[47,6] = 0; $data[47,7] = 9
$data[48,0] = 'Islandia'; $data[48,1] = 1486; $data[48,2] = 69; $data[48,3] = 428; $data[48,4] = 1054; $data[48,5] = 11; $data[48,6] = 0; $data[48,7] = 4
$data[49,0] = 'Colombia'; $data[49,1] = 1485; $data[49,2] = 79; $data[49,3] = 85; $data[49,4] = 1365; $data[49,5] = 50; $data[49,6] = 3; $data[49,7] = 35
$data[50,0] = 'Argentina'; $data[50,1] = 1451; $data[50,2] = 0; $data[50,3] = 280; $data[50,4] = 1127; $data[50,5] = 86; $data[50,6] = 1; $data[50,7] = 44
$data[51,0] = 'Argelia'; $data[51,1] = 1320; $data[51,2] = 69; $data[51,3] = 90; $data[51,4] = 1078; $data[51,5] = 46; $data[51,6] = 22; $data[51,7] = 152
$data[52,0] = 'Singapur'; $data[52,1] = 1309; $data[52,2] = 120; $data[52,3] = 320; $data[52,4] = 983; $data[52,5] = 25; $data[52,6] = 0; $data[52,7] = 6
$data[53,0] = 'Ucrania'; $data[53,1] = 1251; $data[53,2] = 26; $data[53,3] = 25; $data[53,4] = 1194; $data[53,5] = 16; $data[53,6] = 0; $data[53,7] = 32
$data[54,0] = 'Croacia'; $data[54,1] = 1182; $data[54,2] = 56; $data[54,3] = 125; $data[54,4] = 1042; $data[54,5] = 39; $data[54,6] = 3; $data[54,7] = 15
$data[55,0] = 'Egipto'; $data[55,1] = 1173; $data[55,2] = 103; $data[55,3] = 247; $data[55,4] = 848; $data[55,5] = 0; $data[55,6] = 7; $data[55,7] = 78
$data[56,0] = 'Estonia'; $data[56,1] = 1097; $data[56,2] = 58; $data[56,3] = 62; $data[56,4] = 1020; $data[56,5] = 17; $data[56,6] = 2; $data[56,7] = 15
$data[57,0] = 'Nueva Zelanda'; $data[57,1] = 1039; $data[57,2] = 89; $data[57,3] = 156; $data[57,4] = 882; $data[57,5] = 3; $data[57,6] = 0; $data[57,7] = 1
$data[58,0] = 'Eslovenia'; $data[58,1] = 997; $data[58,2] = 20; $data[58,3] = 79; $data[58,4] = 890; $data[58,5] = 31; $data[58,6] = 6; $data[58,7] = 28
$data[59,0] = 'Marruecos'; $data[59,1] = 990; $data[59,2] = 71; $data[59,3] = 71; $data[59,4] = 850; $data[59,5] = 1; $data[59,6] = 10; $data[59,7] = 69
$data[60,0] = 'Irak'; $data[60,1] = 961; $data[60,2] = 83; $data[60,3] = 279; $data[60,4] = 621; $data[60,5] = 0; $data[60,6] = 5; $data[60,7] = 61
$data[61,0] = 'Hong Kong'; $data[61,1] = 890; $data[61,2] = 28; $data[61,3] = 206; $data[61,4] = 680; $data[61,5] = 8; $data[61,6] = 0; $data[61,7] = 4
$data[62,0] = 'Moldavia'; $data[62,1] = 864; $data[62,2] = 112; $data[62,3] = 30; $data[62,4] = 819; $data[62,5] = 80; $data[62,6] = 3; $data[62,7] = 15
$data[63,0] = 'Armenia'; $data[63,1] = 822; $data[63,2] = 52; $data[63,3] = 57; $data[63,4] = 758; $data[63,5] = 30; $data[63,6] = 0; $data[63,7] = 7
$data[64,0] = 'Lituania'; $data[64,1] = 811; $data[64,2] = 40; $data[64,3] = 7; $data[64,4] = 791; $data[64,5] = 11; $data[64,6] = 2; $data[64,7] = 13
$data[65,0] = 'Hungria'; $data[65,1] = 733; $data[65,2] = 55; $data[65,3] = 66; $data[65,4] = 633; $data[65,5] = 17; $data[65,6] = 2; $data[65,7] = 34
$data[66,0] = 'Crucero'; $data[66,1] = 712; $data[66,2] = 0; $data[66,3] = 619; $data[66,4] = 82; $data[66,5] = 10; $data[66,6] = 0; $data[66,7] = 11
$data[67,0] = 'Barein'; $data[67,1] = 698; $data[67,2] = 10; $data[67,3] = 427; $data[67,4] = 267; $data[67,5] = 3; $data[67,6] = 0; $data[67,7] = 4
$data[68,0] = 'Bosnia y Herzegovina'; $data[68,1] = 654; $data[68,2] = 30; $data[68,3] = 30; $data[68,4] = 601; $data[68,5] = 4; $data[68,6] = 2; $data[68,7] = 23
$data[69,0] = 'Azerbaiyan'; $data[69,1] = 584; $data[69,2] = 63; $data[69,3] = 32; $data[69,4] = 545; $data[69,5] = 17; $data[69,6] = 2; $data[69,7] = 7
$data[70,0] = 'Kazajistan'; $data[70,1] = 584; $data[70,2] = 53; $data[70,3] = 42; $data[70,4] = 536; $data[70,5] = 6; $data[70,6] = 1; $data[70,7] = 6
$data[71,0] = 'Tunez'; $data[71,1] = 574; $data[71,2] = 21; $data[71,3] = 5; $data[71,4] = 547; $data[71,5] = 39; $data[71,6] = 4; $data[71,7] = 22
$data[72,0] = 'Bielorrusia'; $data[72,1] = 562; $data[72,2] = 122; $data[72,3] = 52; $data[72,4] = 502; $data[72,5] = 11; $data[72,6] = 3; $data[72,7] = 8
$data[73,0] = 'Kuwait'; $data[73,1] = 556; $data[73,2] = 77; $data[73,3] = 99; $data[73,4] = 456; $data[73,5] = 17; $data[73,6] = 0; $data[73,7] = 1
$data[74,0] = 'Camerun'; $data[74,1] = 555; $data[74,2] = 0; $data[74,3] = 17; $data[74,4] = 529; $data[74,5] = 0; $data[74,6] = 0; $data[74,7] = 9
$data[75,0] = 'Republica de Macedonia'; $data[75,1] = 555; $data[75,2] = 72; $data[75,3] = 23; $data[75,4] = 514; $data[75,5] = 15; $data[75,6] = 1; $data[75,7] = 18
$data[76,0] = 'Letonia'; $data[76,1] = 533; $data[76,2] = 24; $data[76,3] = 1; $data[76,4] = 531; $data[76,5] = 4; $data[76,6] = 0; $data[76,7] = 1
$data[77,0] = 'Bulgaria'; $data[77,1] = 531; $data[77,2] = 28; $data[77,3] = 37; $data[77,4] = 474; $data[77,5] = 22; $data[77,6] = 3; $data[77,7] = 20
$data[78,0] = 'Libano'; $data[78,1] = 527; $data[78,2] = 7; $data[78,3] = 54; $data[78,4] = 455; $data[78,5] = 28; $data[78,6] = 1; $data[78,7] = 18
$data[79,0] = 'Principado de Andorra'; $data[79,1] = 501; $data[79,2] = 35; $data[79,3] = 26; $data[79,4] = 457; $data[79,5] = 12; $data[79,6] = 1; $data[79,7] = 18
$data[80,0] = 'Eslovaquia'; $data[80,1] = 485; $data[80,2] = 14; $data[80,3] = 10; $data[80,4] = 474; $data[80,5] = 3; $data[80,6] = 0; $data[80,7] = 1
$data[81,0] = 'Republica de Chipre'; $data[81,1] = 446; $data[81,2] = 20; $data[81,3] = 37; $data[81,4] = 400; $data[81,5] = 11; $data[81,6] = 0; $data[81,7] = 9
$data[82,0] = 'Costa Rica'; $data[82,1] = 435; $data[82,2] = 0; $data[82,3] = 13; $data[82,4] = 420; $data[82,5] = 13; $data[82,6] = 0; $data[82,7] = 2
$data[83,0] = 'Uruguay'; $data[83,1] = 400; $data[83,2] = 0; $data[83,3] = 93; $data[83,4] = 302; $data[83,5] = 13; $data[83,6] = 0; $data[83,7] = 5
$data[84,0] = 'Taiwan'; $data[84,1] = 363; $data[84,2] = 8; $data[84,3] = 54; $data[84,4] = 304; $data[84,5] = 0; $data[84,6] = 0; $data[84,7] = 5
$data[85,0] = 'Albania'; $data[85,1] = 361; $data[85,2] = 28; $data[85,3] = 104; $data[85,4] = 237; $data[85,5] = 7; $data[85,6] = 0; $data[85,7] = 20
$data[86,0] = 'Afganistan'; $data[86,1] = 349; $data[86,2] = 40; $data[86,3] = 15; $data[86,4] = 327; $data[86,5] = 0; $data[86,6] = 0; $data[86,7] = 7
$data[87,0] = 'Burkina Faso'; $data[87,1] = 345; $data[87,2] = 27; $data[87,3] = 90; $data[87,4] = 238; $data[87,5] = 0; $data[87,6] = 1; $data[87,7] = 17
$data[88,0] = 'Jordania'; $data[88,1] = 345; $data[88,2] = 22; $data[88,3] = 110; $data[88,4] = 230; $data[88,5] = 5; $data[88,6] = 0; $data[88,7] = 5
$data[89,0] = 'Reunion'; $data[89,1] = 344; $data[89,2] = 10; $data[89,3] = 40; $data[89,4] = 304; $data[89,5] = 4; $data[89,6] = 0; $data[89,7] = 0
$data[90,0] = 'Cuba'; $data[90,1] = 320; $data[90,2] = 32; $data[90,3] = 15; $data[90,4] = 297; $data[90,5] = 11; $data[90,6] = 2; $data[90,7] = 8
$data[91,0] = 'Uzbekistan'; $data[91,1] = 310; $data[91,2] = 44; $data[91,3] = 30; $data[91,4] = 278; $data[91,5] = 8; $data[91,6] = 0; $data[91,7] = 2
$data[92,0] = 'Oman'; $data[92,1] = 298; $data[92,2] = 21; $data[92,3] = 61; $data[92,4] = 235; $data[92,5] = 3; $data[92,6] = 0; $data[92,7] = 2
$data[93,0] = 'Honduras'; $data[93,1] = 268; $data[93,2] = 4; $data[93,3] = 6; $data[93,4] = 240; $data[93,5] = 10; $data[93,6] = 7; $data[93,7] = 22
$data[94,0] = 'San Marino'; $data[94,1] = 266; $data[94,2] = 7; $data[94,3] = 35; $data[94,4] = 199; $data[94,5] = 14; $data[94,6] = 0; $data[94,7] = 32
$data[95,0] = 'Costa de Marfil'; $data[95,1] = 245; $data[95,2] = 0; $data[95,3] = 25; $data[95,4] = 219; $data[95,5] = 0; $data[95,6] = 0; $data[95,7] = 1
$data[96,0] = 'Vietnam'; $data[96,1] = 241; $data[96,2] = 1; $data[96,3] = 90; $data[96,4] = 151; $data[96,5] = 3; $data[96,6] = 0; $data[96,7] = 0
$data[97,0] = 'Estado de Palestina'; $data[97,1] = 234; $data[97,2] = 17; $data[97,3] = 23; $data[97,4] = 210; $data[97,5] = 0; $data[97,6] = 0; $data[97,7] = 1
$data[98,0] = 'Malta'; $data[98,1] = 227; $data[98,2] = 14; $data[98,3] = 5; $data[98,4] = 222; $data[98,5] = 3; $data[98,6] = 0; $data[98,7] = 0
$data[99,0] = 'Mauricio'; $data[99,1] = 227; $data[99,2] = 31; $data[99,3] = 7; $data[99,4] = 213; $data[99,5] = 1; $data[99,6] = 0; $data[99,7] = 7
$data[100,0] = 'Nigeria'; $data[100,1] = 224; $data[100,2] = 10; $data[100,3] = 27; $data[100,4] = 192; $data[100,5] = 2; $data[100,6] = 1; $data[100,7] = 5
$data[101,0] = 'Senegal'; $data[101,1] = 222; $data[101,2] = 3; $data[101,3] = 82; $data[101,4] = 138; $data[101,5] = 1; $data[101,6] = 0; $data[101,7] = 2
$data[102,0] = 'Montenegro'; $data[102,1] = 214; $data[102,2] = 13; $data[102,3] = 1; $data[102,4] = 211; $data[102,5] = 4; $data[102,6] = 0; $data[102,7] = 2
$data[103,0] = 'Ghana'; $data[103,1] = 214; $data[103,2] = 9; $data[103,3] = 31; $data[103,4] = 178; $data[103,5] = 2; $data[103,6] = 0; $data[103,7] = 5
$data[104,0] = 'Islas Feroe'; $data[104,1] = 181; $data[104,2] = 0; $data[104,3] = 99; $data[104,4] = 82; $data[104,5] = 1; $data[104,6] = 0; $data[104,7] = 0
$data[105,0] = 'Sri Lanka'; $data[105,1] = 176; $data[105,2] = 10; $data[105,3] = 33; $data[105,4] = 138; $data[105,5] = 5; $data[105,6] = 0; $data[105,7] = 5
$data[106,0] = 'Georgia'; $data[106,1] = 174; $data[106,2] = 12; $data[106,3] = 36; $data[106,4] = 136; $data[106,5] = 6; $data[106,6] = 1; $data[106,7] = 2
$data[107,0] = 'Venezuela'; $data[107,1] = 159; $data[107,2] = 4; $data[107,3] = 52; $data[107,4] = 100; $data[107,5] = 6; $data[107,6] = 0; $data[107,7] = 7
$data[108,0] = 'Bolivia'; $data[108,1] = 157; $data[108,2] = 18; $data[108,3] = 2; $data[108,4] = 145; $data[108,5] = 3; $data[108,6] = 0; $data[108,7] = 10
$data[109,0] = 'Consejo Danes para los Refugiados'; $data[109,1] = 154; $data[109,2] = 0; $data[109,3] = 3; $data[109,4] = 133; $data[109,5] = 0; $data[109,6] = 0; $data[109,7] = 18
$data[110,0] = 'Martinica'; $data[110,1] = 149; $data[110,2] = 4; $data[110,3] = 50; $data[110,4] = 95; $data[110,5] = 21; $data[110,6] = 1; $data[110,7] = 4
$data[111,0] = 'Kirguistan'; $data[111,1] = 147; $data[111,2] = 3; $data[111,3] = 9; $data[111,4] = 137; $data[111,5] = 5; $data[111,6] = 0; $data[111,7] = 1
$data[112,0] = 'Niger'; $data[112,1] = 144; $data[112,2] = 0; $data[112,3] = 0; $data[112,4] = 136; $data[112,5] = 0; $data[112,6] = 0; $data[112,7] = 8
$data[113,0] = 'Kenia'; $data[113,1] = 142; $data[113,2] = 16; $data[113,3] = 4; $data[113,4] = 134; $data[113,5] = 2; $data[113,6] = 0; $data[113,7] = 4
$data[114,0] = 'Brunei'; $data[114,1] = 135; $data[114,2] = 0; $data[114,3] = 73; $data[114,4] = 61; $data[114,5] = 3; $data[114,6] = 0; $data[114,7] = 1
$data[115,0] = 'Mayotte'; $data[115,1] = 134; $data[115,2] = 0; $data[115,3] = 14; $data[115,4] = 118; $data[115,5] = 3; $data[115,6] = 0; $data[115,7] = 2
$data[116,0] = 'Guadalupe'; $data[116,1] = 134; $data[116,2] = 0; $data[116,3] = 24; $data[116,4] = 103; $data[116,5] = 14; $data[116,6] = 0; $data[116,7] = 7
$data[117,0] = 'Isla de Man'; $data[117,1] = 127; $data[117,2] = 1; $data[117,3] = 0; $data[117,4] = 126; $data[117,5] = 0; $data[117,6] = 0; $data[117,7] = 1
$data[118,0] = 'Guinea'; $data[118,1] = 121; $data[118,2] = 10; $data[118,3] = 5; $data[118,4] = 116; $data[118,5] = 0; $data[118,6] = 0; $data[118,7] = 0
$data[119,0] = 'Camboya'; $data[119,1] = 114; $data[119,2] = 0; $data[119,3] = 50; $data[119,4] = 64; $data[119,5] = 1; $data[119,6] = 0; $data[119,7] = 0
$data[120,0] = 'Ruanda'; $data[120,1] = 104; $data[120,2] = 2; $data[120,3] = 4; $data[120,4] = 100; $data[120,5] = 0; $data[120,6] = 0; $data[120,7] = 0
$data[121,0] = 'Trinidad yTobago'; $data[121,1] = 104; $data[121,2] = 1; $data[121,3] = 1; $data[121,4] = 96; $data[121,5] = 0; $data[121,6] = 1; $data[121,7] = 7
$data[122,0] = 'Paraguay'; $data[122,1] = 104; $data[122,2] = 8; $data[122,3] = 12; $data[122,4] = 89; $data[122,5] = 2; $data[122,6] = 0; $data[122,7] = 3
$data[123,0] = 'Gibraltar'; $data[123,1] = 103; $data[123,2] = 5; $data[123,3] = 52; $data[123,4] = 51; $data[123,5] = 0; $data[123,6] = 0; $data[123,7] = 0
$data[124,0] = 'Banglades'; $data[124,1] = 88; $data[124,2] = 18; $data[124,3] = 33; $data[124,4] = 46; $data[124,5] = 1; $data[124,6] = 1; $data[124,7] = 9
$data[125,0] = 'Liechtenstein'; $data[125,1] = 77; $data[125,2] = 0; $data[125,3] = 0; $data[125,4] = 76; $data[125,5] = 0; $data[125,6] = 0; $data[125,7] = 1
$data[126,0] = 'Monaco'; $data[126,1] = 73; $data[126,2] = 7; $data[126,3] = 3; $data[126,4] = 69; $data[126,5] = 2; $data[126,6] = 0; $data[126,7] = 1
$data[127,0] = 'Madagascar'; $data[127,1] = 72; $data[127,2] = 2; $data[127,3] = 2; $data[127,4] = 70; $data[127,5] = 6; $data[127,6] = 0; $data[127,7] = 0
$data[128,0] = 'Aruba'; $data[128,1] = 64; $data[128,2] = 0; $data[128,3] = 1; $data[128,4] = 63; $data[128,5] = 0; $data[128,6] = 0; $data[128,7] = 0
$data[129,0] = 'El Salvador'; $data[129,1] = 62; $data[129,2] = 6; $data[129,3] = 2; $data[129,4] = 57; $data[129,5] = 4; $data[129,6] = 0; $data[129,7] = 3
$data[130,0] = 'Guatemala'; $data[130,1] = 61; $data[130,2] = 0; $data[130,3] = 15; $data[130,4] = 44; $data[130,5] = 1; $data[130,6] = 0; $data[130,7] = 2
$data[131,0] = 'Guayana Francesa'; $data[131,1] = 61; $data[131,2] = 0; $data[131,3] = 22; $data[131,4] = 39; $data[131,5] = 1; $data[131,6] = 0; $data[131,7] = 0
$data[132,0] = 'Republica de Yibuti'; $data[132,1] = 59; $data[132,2] = 9; $data[132,3] = 9; $data[132,4] = 50; $data[132,5] = 0; $data[132,6] = 0; $data[132,7] = 0
$data[133,0] = 'Barbados'; $data[133,1] = 56; $data[133,2] = 4; $data[133,3] = 6; $data[133,4] = 49; $data[133,5] = 4; $data[133,6] = 1; $data[133,7] = 1
$data[134,0] = 'Jamaica'; $data[134,1] = 55; $data[134,2] = 2; $data[134,3] = 7; $data[134,4] = 45; $data[134,5] = 0; $data[134,6] = 0; $data[134,7] = 3
$data[135,0] = 'Uganda'; $data[135,1] = 48; $data[135,2] = 0; $data[135,3] = 0; $data[135,4] = 48; $data[135,5] = 0; $data[135,6] = 0; $data[135,7] = 0
$data[136,0] = 'Mali'; $data[136,1] = 45; $data[136,2] = 4; $data[136,3] = 1; $data[136,4] = 39; $data[136,5] = 0; $data[136,6] = 2; $data[136,7] = 5
$data[137,0] = 'Congo'; $data[137,1] = 45; $data[137,2] = 23; $data[137,3] = 2; $data[137,4] = 38; $data[137,5] = 0; $data[137,6] = 3; $data[137,7] = 5
$data[138,0] = 'Macao'; $data[138,1] = 44; $data[138,2] = 0; $data[138,3] = 10; $data[138,4] = 34; $data[138,5] = 0; $data[138,6] = 0; $data[138,7] = 0
$data[139,0] = 'Togo'; $data[139,1] = 44; $data[139,2] = 3; $data[139,3] = 20; $data[139,4] = 21; $data[139,5] = 0; $data[139,6] = 0; $data[139,7] = 3
$data[140,0] = 'Etiopia'; $data[140,1] = 43; $data[140,2] = 5; $data[140,3] = 4; $data[140,4] = 37; $data[140,5] = 1; $data[140,6] = 2; $data[140,7] = 2
$data[141,0] = 'Polinesia Francesa'; $data[141,1] = 41; $data[141,2] = 1; $data[141,3] = 0; $data[141,4] = 41; $data[141,5] = 1; $data[141,6] = 0; $data[141,7] = 0
$data[142,0] = 'Puerto Rico'; $data[142,1] = 39; $data[142,2] = 0; $data[142,3] = 1; $data[142,4] = 36; $data[142,5] = 0; $data[142,6] = 0; $data[142,7] = 2
$data[143,0] = 'Zambia'; $data[143,1] = 39; $data[143,2] = 0; $data[143,3] = 3; $data[143,4] = 35; $data[143,5] = 0; $data[143,6] = 0; $data[143,7] = 1
$data[144,0] = 'Bermudas'; $data[144,1] = 37; $data[144,2] = 0; $data[144,3] = 14; $data[144,4] = 23; $data[144,5] = 0; $data[144,6] = 0; $data[144,7] = 0
$data[145,0] = 'Islas Caimanes'; $data[145,1] = 35; $data[145,2] = 0; $data[145,3] = 1; $data[145,4] = 33; $data[145,5] = 0; $data[145,6] = 0; $data[145,7] = 1
$data[146,0] = 'Guam'; $data[146,1] = 32; $data[146,2] = 0; $data[146,3] = 0; $data[146,4] = 31; $data[146,5] = 0; $data[146,6] = 0; $data[146,7] = 1
$data[147,0] = 'Eritrea'; $data[147,1] = 29; $data[147,2] = 0; $data[147,3] = 0; $data[147,4] = 29; $data[147,5] = 0; $data[147,6] = 0; $data[147,7] = 0
$data[148,0] = 'San Martin (Parte Francesa)'; $data[148,1] = 29; $data[148,2] = 5; $data[148,3] = 7; $data[148,4] = 20; $data[148,5] = 6; $data[148,6] = 0; $data[148,7] = 2
$data[149,0] = 'Bahamas'; $data[149,1] = 28; $data[149,2] = 0; $data[149,3] = 0; $data[149,4] = 24; $data[149,5] = 1; $data[149,6] = 0; $data[149,7] = 4
$data[150,0] = 'San Martin (Parte Holandesa)'; $data[150,1] = 25; $data[150,2] = 2; $data[150,3] = 6; $data[150,4] = 15; $data[150,5] = 0; $data[150,6] = 0; $data[150,7] = 4
$data[151,0] = 'Guyana'; $data[151,1] = 24; $data[151,2] = 1; $data[151,3] = 0; $data[151,4] = 20; $data[151,5] = 0; $data[151,6] = 0; $data[151,7] = 4
$data[152,0] = 'Benin'; $data[152,1] = 22; $data[152,2] = 6; $data[152,3] = 5; $data[152,4] = 17; $data[152,5] = 0; $data[152,6] = 0; $data[152,7] = 0
$data[153,0] = 'Birmania'; $data[153,1] = 21; $data[153,2] = 0; $data[153,3] = 0; $data[153,4] = 20; $data[153,5] = 0; $data[153,6] = 0; $data[153,7] = 1
$data[154,0] = 'Gabon'; $data[154,1] = 21; $data[154,2] = 0; $data[154,3] = 1; $data[154,4] = 19; $data[154,5] = 0; $data[154,6] = 0; $data[154,7] = 1
$data[155,0] = 'Haiti'; $data[155,1] = 21; $data[155,2] = 1; $data[155,3] = 1; $data[155,4] = 19; $data[155,5] = 0; $data[155,6] = 1; $data[155,7] = 1
$data[156,0] = 'Tanzania'; $data[156,1] = 20; $data[156,2] = 0; $data[156,3] = 3; $data[156,4] = 16; $data[156,5] = 0; $data[156,6] = 0; $data[156,7] = 1
$data[157,0] = 'Siria'; $data[157,1] = 19; $data[157,2] = 3; $data[157,3] = 2; $data[157,4] = 15; $data[157,5] = 0; $data[157,6] = 0; $data[157,7] = 2
$data[158,0] = 'Maldivas'; $data[158,1] = 19; $data[158,2] = 0; $data[158,3] = 13; $data[158,4] = 6; $data[158,5] = 0; $data[158,6] = 0; $data[158,7] = 0
$data[159,0] = 'Guinea-Bisau'; $data[159,1] = 18; $data[159,2] = 0; $data[159,3] = 0; $data[159,4] = 18; $data[159,5] = 0; $data[159,6] = 0; $data[159,7] = 0
$data[160,0] = 'Libia'; $data[160,1] = 18; $data[160,2] = 0; $data[160,3] = 0; $data[160,4] = 17; $data[160,5] = 0; $data[160,6] = 0; $data[160,7] = 1
$data[161,0] = 'Islas Virgenes de los Estados Unidos'; $data[161,1] = 17; $data[161,2] = 0; $data[161,3] = 0; $data[161,4] = 17; $data[161,5] = 0; $data[161,6] = 0; $data[161,7] = 0
$data[162,0] = 'Nueva Caledonia'; $data[162,1] = 17; $data[162,2] = 0; $data[162,3] = 1; $data[162,4] = 16; $data[162,5] = 0; $data[162,6] = 0; $data[162,7] = 0
$data[163,0] = 'Guinea Ecuatorial'; $data[163,1] = 16; $data[163,2] = 0; $data[163,3] = 1; $data[163,4] = 15; $data[163,5] = 0; $data[163,6] = 0; $data[163,7] = 0
$data[164,0] = 'Namibia'; $data[164,1] = 16; $data[164,2] = 2; $data[164,3] = 3; $data[164,4] = 13; $data[164,5] = 0; $data[164,6] = 0; $data[164,7] = 0
$data[165,0] = 'Antigua y Barbuda'; $data[165,1] = 15; $data[165,2] = 0; $data[165,3] = 0; $data[165,4] = 15; $data[165,5] = 1; $data[165,6] = 0; $data[165,7] = 0
$data[166,0] = 'Dominica'; $data[166,1] = 14; $data[166,2] = 0; $data[166,3] = 0; $data[166,4] = 14; $data[166,5] = 0; $data[166,6] = 0; $data[166,7] = 0
$data[167,0] = 'Santa Lucia'; $data[167,1] = 14; $data[167,2] = 0; $data[167,3] = 1; $data[167,4] = 13; $data[167,5] = 0; $data[167,6] = 0; $data[167,7] = 0
$data[168,0] = 'Mongolia'; $data[168,1] = 14; $data[168,2] = 0; $data[168,3] = 2; $data[168,4] = 12; $data[168,5] = 0; $data[168,6] = 0; $data[168,7] = 0
$data[169,0] = 'Liberia'; $data[169,1] = 13; $data[169,2] = 3; $data[169,3] = 3; $data[169,4] = 7; $data[169,5] = 0; $data[169,6] = 2; $data[169,7] = 3
$data[170,0] = 'Fiyi'; $data[170,1] = 12; $data[170,2] = 0; $data[170,3] = 0; $data[170,4] = 12; $data[170,5] = 0; $data[170,6] = 0; $data[170,7] = 0
$data[171,0] = 'Granada'; $data[171,1] = 12; $data[171,2] = 0; $data[171,3] = 0; $data[171,4] = 12; $data[171,5] = 2; $data[171,6] = 0; $data[171,7] = 0
$data[172,0] = 'Sudan'; $data[172,1] = 12; $data[172,2] = 2; $data[172,3] = 2; $data[172,4] = 8; $data[172,5] = 0; $data[172,6] = 0; $data[172,7] = 2
$data[173,0] = 'Laos'; $data[173,1] = 11; $data[173,2] = 1; $data[173,3] = 0; $data[173,4] = 11; $data[173,5] = 0; $data[173,6] = 0; $data[173,7] = 0
$data[174,0] = 'Groenlandia'; $data[174,1] = 11; $data[174,2] = 0; $data[174,3] = 3; $data[174,4] = 8; $data[174,5] = 0; $data[174,6] = 0; $data[174,7] = 0
$data[175,0] = 'Curazao'; $data[175,1] = 11; $data[175,2] = 0; $data[175,3] = 5; $data[175,4] = 5; $data[175,5] = 0; $data[175,6] = 0; $data[175,7] = 1
$data[176,0] = 'San Cristobal y Nieves'; $data[176,1] = 10; $data[176,2] = 1; $data[176,3] = 0; $data[176,4] = 10; $data[176,5] = 0; $data[176,6] = 0; $data[176,7] = 0
$data[177,0] = 'Seychelles'; $data[177,1] = 10; $data[177,2] = 0; $data[177,3] = 0; $data[177,4] = 10; $data[177,5] = 0; $data[177,6] = 0; $data[177,7] = 0
$data[178,0] = 'Surinam'; $data[178,1] = 10; $data[178,2] = 0; $data[178,3] = 0; $data[178,4] = 9; $data[178,5] = 0; $data[178,6] = 0; $data[178,7] = 1
$data[179,0] = 'Mozambique'; $data[179,1] = 10; $data[179,2] = 0; $data[179,3] = 1; $data[179,4] = 9; $data[179,5] = 0; $data[179,6] = 0; $data[179,7] = 0
$data[180,0] = 'Angola'; $data[180,1] = 10; $data[180,2] = 0; $data[180,3] = 2; $data[180,4] = 6; $data[180,5] = 0; $data[180,6] = 0; $data[180,7] = 2
$data[181,0] = 'Republica del Chad'; $data[181,1] = 9; $data[181,2] = 0; $data[181,3] = 0; $data[181,4] = 9; $data[181,5] = 0; $data[181,6] = 0; $data[181,7] = 0
$data[182,0] = 'Suazilandia'; $data[182,1] = 9; $data[182,2] = 0; $data[182,3] = 0; $data[182,4] = 9; $data[182,5] = 0; $data[182,6] = 0; $data[182,7] = 0
$data[183,0] = 'Nepal'; $data[183,1] = 9; $data[183,2] = 0; $data[183,3] = 1; $data[183,4] = 8; $data[183,5] = 0; $data[183,6] = 0; $data[183,7] = 0
$data[184,0] = 'Zimbabue'; $data[184,1] = 9; $data[184,2] = 0; $data[184,3] = 0; $data[184,4] = 8; $data[184,5] = 0; $data[184,6] = 0; $data[184,7] = 1
$data[185,0] = 'Montserrat'; $data[185,1] = 9; $data[185,2] = 0; $data[185,3] = 0; $data[185,4] = 7; $data[185,5] = 0; $data[185,6] = 0; $data[185,7] = 2
$data[186,0] = 'Republica de Africa Central'; $data[186,1] = 8; $data[186,2] = 0; $data[186,3] = 0; $data[186,4] = 8; $data[186,5] = 0; $data[186,6] = 0; $data[186,7] = 0
$data[187,0] = 'Santa Sede'; $data[187,1] = 7; $data[187,2] = 0; $data[187,3] = 0; $data[187,4] = 7; $data[187,5] = 0; $data[187,6] = 0; $data[187,7] = 0
$data[188,0] = 'San Vicente y las Granadinas'; $data[188,1] = 7; $data[188,2] = 0; $data[188,3] = 1; $data[188,4] = 6; $data[188,5] = 0; $data[188,6] = 0; $data[188,7] = 0
$data[189,0] = 'Somalia'; $data[189,1] = 7; $data[189,2] = 0; $data[189,3] = 1; $data[189,4] = 6; $data[189,5] = 0; $data[189,6] = 0; $data[189,7] = 0
$data[190,0] = 'Cabo Verde'; $data[190,1] = 7; $data[190,2] = 0; $data[190,3] = 0; $data[190,4] = 6; $data[190,5] = 0; $data[190,6] = 0; $data[190,7] = 1
$data[191,0] = 'Sierra Leona'; $data[191,1] = 6; $data[191,2] = 2; $data[191,3] = 0; $data[191,4] = 6; $data[191,5] = 0; $data[191,6] = 0; $data[191,7] = 0
$data[192,0] = 'Botsuana'; $data[192,1] = 6; $data[192,2] = 2; $data[192,3] = 0; $data[192,4] = 5; $data[192,5] = 0; $data[192,6] = 0; $data[192,7] = 1
$data[193,0] = 'San Bartolome'; $data[193,1] = 6; $data[193,2] = 0; $data[193,3] = 1; $data[193,4] = 5; $data[193,5] = 0; $data[193,6] = 0; $data[193,7] = 0
$data[194,0] = 'Nicaragua'; $data[194,1] = 6; $data[194,2] = 1; $data[194,3] = 0; $data[194,4] = 5; $data[194,5] = 0; $data[194,6] = 0; $data[194,7] = 1
$data[195,0] = 'Mauritania'; $data[195,1] = 6; $data[195,2] = 0; $data[195,3] = 2; $data[195,4] = 3; $data[195,5] = 0; $data[195,6] = 0; $data[195,7] = 1
$data[196,0] = 'Belice'; $data[196,1] = 5; $data[196,2] = 1; $data[196,3] = 0; $data[196,4] = 5; $data[196,5] = 1; $data[196,6] = 0; $data[196,7] = 0
$data[197,0] = 'Islas Turcas y Caicos'; $data[197,1] = 5; $data[197,2] = 0; $data[197,3] = 0; $data[197,4] = 4; $data[197,5] = 0; $data[197,6] = 1; $data[197,7] = 1
$data[198,0] = 'Butan'; $data[198,1] = 5; $data[198,2] = 0; $data[198,3] = 2; $data[198,4] = 3; $data[198,5] = 0; $data[198,6] = 0; $data[198,7] = 0
$data[199,0] = 'Sahara Occidental'; $data[199,1] = 4; $data[199,2] = 0; $data[199,3] = 0; $data[199,4] = 4; $data[199,5] = 0; $data[199,6] = 0; $data[199,7] = 0
$data[200,0] = 'Malaui'; $data[200,1] = 4; $data[200,2] = 0; $data[200,3] = 0; $data[200,4] = 4; $data[200,5] = 0; $data[200,6] = 0; $data[200,7] = 0
$data[201,0] = 'Gambia'; $data[201,1] = 4; $data[201,2] = 0; $data[201,3] = 2; $data[201,4] = 1; $data[201,5] = 0; $data[201,6] = 0; $data[201,7] = 1
$data[202,0] = 'Islas Virgenes Britanicas'; $data[202,1] = 3; $data[202,2] = 0; $data[202,3] = 0; $data[202,4] = 3; $data[202,5] = 0; $data[202,6] = 0; $data[202,7] = 0
$data[203,0] = 'Anguila'; $data[203,1] = 3; $data[203,2] = 0; $data[203,3] = 0; $data[203,4] = 3; $data[203,5] = 0; $data[203,6] = 0; $data[203,7] = 0
$data[204,0] = 'Burundi'; $data[204,1] = 3; $data[204,2] = 0; $data[204,3] = 0; $data[204,4] = 3; $data[204,5] = 0; $data[204,6] = 0; $data[204,7] = 0
$data[205,0] = 'Islas Malvinas'; $data[205,1] = 2; $data[205,2] = 1; $data[205,3] = 0; $data[205,4] = 2; $data[205,5] = 0; $data[205,6] = 0; $data[205,7] = 0
$data[206,0] = 'Bonaire, San Eustaquio y Saba'; $data[206,1] = 2; $data[206,2] = 0; $data[206,3] = 0; $data[206,4] = 2; $data[206,5] = 0; $data[206,6] = 0; $data[206,7] = 0
$data[207,0] = 'Sudan del Sur'; $data[207,1] = 1; $data[207,2] = 1; $data[207,3] = 0; $data[207,4] = 1; $data[207,5] = 0; $data[207,6] = 0; $data[207,7] = 0
$data[208,0] = 'Timor Oriental'; $data[208,1] = 1; $data[208,2] = 0; $data[208,3] = 0; $data[208,4] = 1; $data[208,5] = 0; $data[208,6] = 0; $data[208,7] = 0
$data[209,0] = 'Papua Nueva Guinea'; $data[209,1] = 1; $data[209,2] = 0; $data[209,3] = 0; $data[209,4] = 1; $data[209,5] = 0; $data[209,6] = 0; $data[209,7] = 0
$data[210,0] = 'San Pedro y Miquelon'; $data[210,1] = 1; $data[210,2] = 1; $data[210,3] = 0; $data[210,4] = 1; $data[210,5] = 0; $data[210,6] = 0; $data[210,7] = 0

$ws.Range("A4:H214").Value = $data

